$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30: predidx / pred_name update
$ws.Range("D30").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E30").Value = "['Normal', 'SoftwareFault']"

# Row 88: predidx / pred_name update
$ws.Range("D88").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E88").Value = "['Normal']"

# Row 113: predidx / pred_name update
$ws.Range("D113").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E113").Value = "['Normal']"
